# Insert a new weekly price record as row 179 (Feria Lagunitas de Puerto Montt - Pomelo).
# Existing rows 179..243 shift down to 180..244; the sheet's used range grows
# from A1:T243 to A1:T244.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(179).Insert()

$ws.Cells.Item(179, 1).Value = 4
$ws.Cells.Item(179, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(179, 3).Value = "Los Lagos"
$ws.Cells.Item(179, 4).Value = 44627
$ws.Cells.Item(179, 5).Value = 10
$ws.Cells.Item(179, 6).Value = "Fruta"
$ws.Cells.Item(179, 7).Value = 100102
$ws.Cells.Item(179, 8).Value = "Cítricos"
$ws.Cells.Item(179, 9).Value = 100102006
$ws.Cells.Item(179, 10).Value = "Pomelo"
$ws.Cells.Item(179, 11).Value = "Start Ruby"
$ws.Cells.Item(179, 12).Value = "Especial"
$ws.Cells.Item(179, 13).Value = 80
$ws.Cells.Item(179, 14).Value = 14000
$ws.Cells.Item(179, 15).Value = 14000
$ws.Cells.Item(179, 16).Value = 14000
$ws.Cells.Item(179, 17).Value = "`$/caja 14 kilos empedrada"
$ws.Cells.Item(179, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(179, 19).Value = 1000
$ws.Cells.Item(179, 20).Value = 14
